$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.942.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5088"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06391"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5463"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7869"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.020.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.426"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.979"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.046"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.873"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.896"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05030"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.267"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.370"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8962"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.622"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5512"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  +14.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.007"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.547"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.648"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8170"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.780.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4537"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05081"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
